$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.041.12"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "2.209.82"
$ws.Range("E3").Value = "  -0.66%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "293.73"
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "87.42"
$ws.Range("E6").Value = "  +1.70%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.509"
$ws.Range("E7").Value = "  -0.99%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.472"
$ws.Range("E9").Value = "  +0.33%  "
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0778"
$ws.Range("E10").Value = "  -1.52%  "
$ws.Range("B11").Value = "Avalanche"
$ws.Range("C11").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "29.90"
$ws.Range("E11").Value = "  -3.63%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "49.73"
$ws.Range("E12").Value = "  +5.81%  "
$ws.Range("E13").Value = "  +2.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.44"
$ws.Range("E14").Value = "  +0.38%  "
$ws.Range("D15").Value = "2.549.72"
$ws.Range("E15").Value = "  -0.64%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.74"
$ws.Range("E16").Value = "  -2.11%  "
$ws.Range("D17").Value = "2.176.74"
$ws.Range("E17").Value = "  -1.99%  "
$ws.Range("E18").Value = "  -0.51%  "
$ws.Range("D19").Value = "39.942.67"
$ws.Range("E19").Value = "  +0.32%  "
$ws.Range("D20").Value = "0.0₃0884"
$ws.Range("E20").Value = "  -0.56%  "
$ws.Range("E21").Value = "  +4.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.76"
$ws.Range("E22").Value = "  -0.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.27"
$ws.Range("E23").Value = "  +0.42%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "238.06"
$ws.Range("E24").Value = "  +1.16%  "
$ws.Range("E25").Value = "  +0.18%  "
$ws.Range("E26").Value = "  -0.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.79"
$ws.Range("E27").Value = "  -2.67%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.46"
$ws.Range("E29").Value = "  -3.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.17"
$ws.Range("E30").Value = "  -0.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "156.85"
$ws.Range("E31").Value = "  +3.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.36"
$ws.Range("E32").Value = "  -5.36%  "
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.89"
$ws.Range("E34").Value = "  +0.61%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0709"
$ws.Range("E35").Value = "  -1.01%  "
$ws.Range("E36").Value = "  -2.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.82"
$ws.Range("E37").Value = "  +3.94%  "
$ws.Range("E38").Value = "  +0.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0978"
$ws.Range("E39").Value = "  -1.65%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "15.34"
$ws.Range("E40").Value = "  -5.91%  "
$ws.Range("E41").Value = "  -1.32%  "
$ws.Range("D42").Value = "2.123.74"
$ws.Range("E42").Value = "  +4.49%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.74"
$ws.Range("E43").Value = "  -2.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.11"
$ws.Range("E44").Value = "  -2.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0267"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.54"
$ws.Range("E46").Value = "  +7.84%  "
$ws.Range("E47").Value = "  -4.26%  "
$ws.Range("E48").Value = "  +3.40%  "
$ws.Range("D49").Value = "2.419.35"
$ws.Range("E49").Value = "  -0.80%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.47"
$ws.Range("E50").Value = "  +2.11%  "
$ws.Range("E51").Value = "  +0.87%  "
